$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) and "全部类型" (All types) both contain the same
# data table; update the "想去人数" (want-to-go count) column F for both.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 51
    $ws.Range("F3").Value = 285
    $ws.Range("F5").Value = 76
}
